# Re-sort / reshuffle the weekly Pepino dulce price rows (rows 2-40)
# so each row reflects the updated date/quality/volume/price data per the
# latest daily-logic consolidation pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = (Get-Date -Year 2021 -Month 4 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 10500
$ws.Range("P2").Value = 583

$ws.Range("D3").Value = (Get-Date -Year 2021 -Month 4 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 8000
$ws.Range("P3").Value = 444

$ws.Range("D4").Value = (Get-Date -Year 2021 -Month 6 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("P4").Value = 750

$ws.Range("D5").Value = (Get-Date -Year 2021 -Month 6 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("P5").Value = 611

$ws.Range("D6").Value = (Get-Date -Year 2021 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 9500
$ws.Range("P6").Value = 528

$ws.Range("D7").Value = (Get-Date -Year 2021 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 8000
$ws.Range("P7").Value = 444

$ws.Range("D8").Value = (Get-Date -Year 2021 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("P8").Value = 861

$ws.Range("D9").Value = (Get-Date -Year 2021 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 13000
$ws.Range("P9").Value = 722

$ws.Range("D10").Value = (Get-Date -Year 2021 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 10500
$ws.Range("P10").Value = 583

$ws.Range("D11").Value = (Get-Date -Year 2021 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 9000
$ws.Range("M11").Value = 9000
$ws.Range("P11").Value = 500

$ws.Range("D12").Value = (Get-Date -Year 2021 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("P12").Value = 861

$ws.Range("D13").Value = (Get-Date -Year 2021 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 13000
$ws.Range("P13").Value = 722

$ws.Range("D14").Value = (Get-Date -Year 2021 -Month 7 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15500
$ws.Range("P14").Value = 861

$ws.Range("D15").Value = (Get-Date -Year 2021 -Month 7 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 14000
$ws.Range("M15").Value = 14000
$ws.Range("P15").Value = 778

$ws.Range("D16").Value = (Get-Date -Year 2021 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 9000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 9500
$ws.Range("P16").Value = 528

$ws.Range("D17").Value = (Get-Date -Year 2021 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 8000
$ws.Range("P17").Value = 444

$ws.Range("D18").Value = (Get-Date -Year 2021 -Month 2 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("P18").Value = 667

$ws.Range("D19").Value = (Get-Date -Year 2021 -Month 2 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 10000
$ws.Range("P19").Value = 556

$ws.Range("D20").Value = (Get-Date -Year 2021 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13500
$ws.Range("P20").Value = 750

$ws.Range("D21").Value = (Get-Date -Year 2021 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 11000
$ws.Range("L21").Value = 11000
$ws.Range("M21").Value = 11000
$ws.Range("P21").Value = 611

$ws.Range("D22").Value = (Get-Date -Year 2022 -Month 2 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 15500
$ws.Range("P22").Value = 861

$ws.Range("D23").Value = (Get-Date -Year 2021 -Month 4 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10500
$ws.Range("P23").Value = 583

$ws.Range("D24").Value = (Get-Date -Year 2021 -Month 4 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I24").Value = "Segunda"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 9000
$ws.Range("P24").Value = 500

$ws.Range("D25").Value = (Get-Date -Year 2021 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("P25").Value = 806

$ws.Range("D26").Value = (Get-Date -Year 2021 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 11000
$ws.Range("M26").Value = 10500
$ws.Range("P26").Value = 583

$ws.Range("D27").Value = (Get-Date -Year 2021 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = 9000
$ws.Range("P27").Value = 500

$ws.Range("D28").Value = (Get-Date -Year 2021 -Month 3 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 11000
$ws.Range("M28").Value = 10500
$ws.Range("P28").Value = 583

$ws.Range("D29").Value = (Get-Date -Year 2021 -Month 3 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I29").Value = "Segunda"
$ws.Range("J29").Value = 50
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = 9000
$ws.Range("P29").Value = 500

$ws.Range("D30").Value = (Get-Date -Year 2021 -Month 3 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 11000
$ws.Range("M30").Value = 10500
$ws.Range("P30").Value = 583

$ws.Range("D31").Value = (Get-Date -Year 2021 -Month 3 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I31").Value = "Segunda"
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 9000
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = 9000
$ws.Range("P31").Value = 500

$ws.Range("D32").Value = (Get-Date -Year 2021 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13500
$ws.Range("P32").Value = 750

$ws.Range("D33").Value = (Get-Date -Year 2021 -Month 5 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 11500
$ws.Range("P33").Value = 639

$ws.Range("D34").Value = (Get-Date -Year 2021 -Month 5 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 50
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = 9000
$ws.Range("P34").Value = 500

$ws.Range("D35").Value = (Get-Date -Year 2021 -Month 3 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 12000
$ws.Range("L35").Value = 13000
$ws.Range("M35").Value = 12500
$ws.Range("P35").Value = 694

$ws.Range("D36").Value = (Get-Date -Year 2021 -Month 3 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I36").Value = "Segunda"
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = 10000
$ws.Range("P36").Value = 556

$ws.Range("D37").Value = (Get-Date -Year 2021 -Month 5 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 12000
$ws.Range("L37").Value = 13000
$ws.Range("M37").Value = 12500
$ws.Range("P37").Value = 694

$ws.Range("D38").Value = (Get-Date -Year 2021 -Month 5 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = 10000
$ws.Range("P38").Value = 556

$ws.Range("D39").Value = (Get-Date -Year 2021 -Month 6 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 11000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = 11500
$ws.Range("P39").Value = 639

$ws.Range("D40").Value = (Get-Date -Year 2021 -Month 6 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I40").Value = "Segunda"
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = 10000
$ws.Range("P40").Value = 556
